$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road | Potion
$ws.Range("H17").Value = 3564.3635
$ws.Range("J17").Value = 3564.3635
$ws.Range("L17").Value = 10693.0905
$ws.Range("N17").Value = -11029.0905

# Row 64: Forged from the Void | Void Glue
$ws.Range("H64").Value = 3556.125
$ws.Range("I64").Value = 3397.8
$ws.Range("J64").Value = 3820
$ws.Range("K64").Value = 3397.8
$ws.Range("L64").Value = 3820
$ws.Range("M64").Value = -3149.8
$ws.Range("N64").Value = -4316

# Row 67: Dodging the Draft (L) | Void Glue
$ws.Range("H67").Value = 3556.125
$ws.Range("I67").Value = 3397.8
$ws.Range("J67").Value = 3820
$ws.Range("K67").Value = 3397.8
$ws.Range("L67").Value = 3820
$ws.Range("M67").Value = -2539.8
$ws.Range("N67").Value = -5536

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 1507.6571
$ws.Range("I137").Value = 1038.421
$ws.Range("J137").Value = 2064.875
$ws.Range("K137").Value = 3115.263
$ws.Range("L137").Value = 6194.625
$ws.Range("M137").Value = -565.2629999999999
$ws.Range("N137").Value = -11294.625

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain''t Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 5710.1
$ws.Range("I2").Value = 6177.8887
$ws.Range("K2").Value = 6177.8887
$ws.Range("M2").Value = -6064.8887

# Row 9: Headbangers'' Thrall | Bronze Barbut
$ws.Range("H9").Value = 8888
$ws.Range("J9").Value = 8888
$ws.Range("L9").Value = 8888
$ws.Range("N9").Value = -9228

# Row 20: Cover Girl | Decorated Bronze Barbut
$ws.Range("H20").Value = 8888
$ws.Range("J20").Value = 8888
$ws.Range("L20").Value = 8888
$ws.Range("N20").Value = -9428

# Row 60: Booty Call | Cobalt-plated Jackboots
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Range("H63").Value = 1887.1282
$ws.Range("I63").Value = 1665.4482
$ws.Range("K63").Value = 1665.4482
$ws.Range("M63").Value = -979.4482

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Range("H66").Value = 1887.1282
$ws.Range("I66").Value = 1665.4482
$ws.Range("K66").Value = 8327.241
$ws.Range("M66").Value = -4895.241

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 5710.1
$ws.Range("I116").Value = 6177.8887
$ws.Range("K116").Value = 6177.8887
$ws.Range("M116").Value = -3883.8887

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 1976.5161
$ws.Range("I122").Value = 1909.7273
$ws.Range("K122").Value = 5729.1819
$ws.Range("M122").Value = -3279.1819

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 5710.1
$ws.Range("I3").Value = 6177.8887
$ws.Range("K3").Value = 6177.8887
$ws.Range("M3").Value = -6063.8887

# Row 5: Axe Me Anything | Bronze War Axe
$ws.Range("H5").Value = 402
$ws.Range("I5").Value = 369.33334
$ws.Range("K5").Value = 369.33334
$ws.Range("M5").Value = -256.33334

# Row 31: When Rhalgr Met Nophica | Spiked Knuckles
$ws.Range("H31").Value = 980
$ws.Range("I31").Value = 980
$ws.Range("K31").Value = 980
$ws.Range("M31").Value = -728

# Row 34: Cleaving the Glim | Iron Round Knife
$ws.Range("H34").Value = 8000
$ws.Range("J34").Value = 8000
$ws.Range("L34").Value = 8000
$ws.Range("N34").Value = -8228

# Row 36: I Saw What You Did There | Iron Chocobotail Saw
$ws.Range("H36").Value = 1137
$ws.Range("I36").Value = 1137
$ws.Range("K36").Value = 1137
$ws.Range("M36").Value = -603

# Row 37: That''s Some Fine Grinding | Initiate''s Mortar
$ws.Range("H37").Value = 5250
$ws.Range("I37").Value = 500
$ws.Range("K37").Value = 500
$ws.Range("M37").Value = -363

# Row 43: Don''t Fear the Reaper | Steel Scythe
$ws.Range("H43").Value = 145000
$ws.Range("J43").Value = 145000
$ws.Range("L43").Value = 145000
$ws.Range("N43").Value = -145362

# Row 87: Winter Weather Conditions | Adamantite Dolabra
$ws.Range("H87").Value = 47000
$ws.Range("J87").Value = 47000
$ws.Range("L87").Value = 47000
$ws.Range("N87").Value = -49496

# Row 90: The Nightsoil Is Dark and Full of Terrors (L) | Adamantite Dolabra
$ws.Range("H90").Value = 47000
$ws.Range("J90").Value = 47000
$ws.Range("L90").Value = 141000
$ws.Range("N90").Value = -153480

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 10869834
$ws.Range("I94").Value = 15625221
$ws.Range("K94").Value = 15625221
$ws.Range("M94").Value = -15624770

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 55556560
$ws.Range("I99").Value = 62500936
$ws.Range("K99").Value = 62500936
$ws.Range("M99").Value = -62499438

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 1333.375
$ws.Range("I107").Value = 986.8570999999999
$ws.Range("J107").Value = 1602.8889
$ws.Range("K107").Value = 986.8570999999999
$ws.Range("L107").Value = 1602.8889
$ws.Range("M107").Value = 933.1429000000001
$ws.Range("N107").Value = -5442.8889

$ws = $wb.Worksheets.Item("CRP")
# Row 20: Re-crating the Scene | Iron Spear
$ws.Range("H20").Value = 46635.7
$ws.Range("J20").Value = 46635.7
$ws.Range("L20").Value = 46635.7
$ws.Range("N20").Value = -47107.7

# Row 30: Polearms Aplenty | Iron Spear
$ws.Range("H30").Value = 46635.7
$ws.Range("J30").Value = 46635.7
$ws.Range("L30").Value = 46635.7
$ws.Range("N30").Value = -46817.7

# Row 92: Walk the Walk | Beech Rod
$ws.Range("H92").Value = 35100.168
$ws.Range("J92").Value = 35100.168
$ws.Range("L92").Value = 35100.168
$ws.Range("N92").Value = -40092.168

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 1647.8
$ws.Range("I99").Value = 1616.2858
$ws.Range("J99").Value = 1721.3334
$ws.Range("K99").Value = 1616.2858
$ws.Range("L99").Value = 1721.3334
$ws.Range("M99").Value = -118.2858000000001
$ws.Range("N99").Value = -4717.3334

# Row 114: Ground to a Halt | White Ash Grinding Wheel
$ws.Range("H114").Value = 31995
$ws.Range("J114").Value = 31995
$ws.Range("L114").Value = 31995
$ws.Range("N114").Value = -40673

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 1647.8
$ws.Range("I126").Value = 1616.2858
$ws.Range("J126").Value = 1721.3334
$ws.Range("K126").Value = 4848.857400000001
$ws.Range("L126").Value = 5164.0002
$ws.Range("M126").Value = -2378.857400000001
$ws.Range("N126").Value = -10104.0002

# Row 128: An A-prop-riate Request | Ironwood Spear
$ws.Range("H128").Value = 46635.7
$ws.Range("J128").Value = 46635.7
$ws.Range("L128").Value = 46635.7
$ws.Range("N128").Value = -56595.7

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 949.05884
$ws.Range("J5").Value = 641.4286
$ws.Range("L5").Value = 1924.2858
$ws.Range("N5").Value = -2148.2858

# Row 45: Don''t Turn Up Your Nose | Sauerkraut
$ws.Range("H45").Value = 689
$ws.Range("J45").Value = 689
$ws.Range("L45").Value = 2067
$ws.Range("N45").Value = -3131

# Row 88: Don''t Let It Fall Apart | Liver-cheese Sandwich
$ws.Range("H88").Value = 6228.5713
$ws.Range("J88").Value = 6228.5713
$ws.Range("L88").Value = 18685.7139
$ws.Range("N88").Value = -19541.7139

# Row 91: Better Come Back with a Sandwich (L) | Liver-cheese Sandwich
$ws.Range("H91").Value = 6228.5713
$ws.Range("J91").Value = 6228.5713
$ws.Range("L91").Value = 18685.7139
$ws.Range("N91").Value = -21649.7139

# Row 97: The Frier Never Lies | Cottonseed Oil
$ws.Range("H97").Value = 683.9
$ws.Range("I97").Value = 756.6667
$ws.Range("J97").Value = 652.7143
$ws.Range("K97").Value = 2270.0001
$ws.Range("L97").Value = 1958.1429
$ws.Range("M97").Value = -1774.0001
$ws.Range("N97").Value = -2950.1429

# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 949.05884
$ws.Range("J135").Value = 641.4286
$ws.Range("L135").Value = 5772.8574
$ws.Range("N135").Value = -10842.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 22: Bad to the Bone | Brass Circlet (Sunstone)
$ws.Range("H22").Value = 8
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = 8
$ws.Range("M22").Value = 521

# Row 109: You''re My Wonderhall | Hematite Earrings of Healing
$ws.Range("H109").Value = 40142.5
$ws.Range("J109").Value = 40142.5
$ws.Range("L109").Value = 40142.5
$ws.Range("N109").Value = -42222.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 1906.8572
$ws.Range("I7").Value = 1832.5834
$ws.Range("K7").Value = 1832.5834
$ws.Range("M7").Value = -1720.5834

# Row 101: A Stitch in Time | Marid Leather Gloves of Healing
$ws.Range("H101").Value = 11665.25
$ws.Range("J101").Value = 11665.25
$ws.Range("L101").Value = 11665.25
$ws.Range("N101").Value = -18155.25

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 1906.8572
$ws.Range("I126").Value = 1832.5834
$ws.Range("K126").Value = 5497.7502
$ws.Range("M126").Value = -3027.7502

$ws = $wb.Worksheets.Item("WVR")
# Row 103: To the Tops | Serge Gambison of Healing
$ws.Range("H103").Value = 19034
$ws.Range("J103").Value = 19034
$ws.Range("L103").Value = 19034
$ws.Range("N103").Value = -21378

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 90912776
$ws.Range("I126").Value = 100003060
$ws.Range("K126").Value = 300009180
$ws.Range("M126").Value = -300006710
